# Reverted default group numbers and metadata so they are consistent with
# the other default paths and data.
# Updates the "UnitMass" values in column C of the "+ loading" and
# "- loading" tables on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "C2"  = 23
    "C3"  = 71
    "C4"  = 39
    "C5"  = 102
    "C6"  = 55
    "C7"  = 46
    "C8"  = 27
    "C9"  = 28
    "C10" = 41
    "C11" = 43
    "C12" = 1
    "C13" = 149
    "C14" = 74
    "C15" = 29
    "C16" = 24
    "C17" = 13
    "C18" = 12
    "C19" = 7
    "C20" = 50
    "C21" = 26
    "C23" = 57
    "C24" = 113
    "C25" = 91
    "C26" = 40
    "C27" = 69
    "C28" = 82
    "C29" = 123
    "C30" = 45
    "C31" = 15
    "C32" = 19
    "C33" = 83
    "C34" = 96
    "C35" = 125
    "C36" = 99
    "C37" = 98
    "C38" = 117
    "C39" = 111
    "C40" = 68
    "C41" = 81
    "C42" = 108
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
